$d = $word.ActiveDocument
$t = $d.Tables(1)

# The worksheet table holds five "problem rows" (1, 5, 10, 15, 20), each with
# five cells of "AxB=C" text that must be replaced with new values. A couple
# of the old values ("39x70=2730") are duplicated across rows, so replacements
# are targeted cell-by-cell via the table's Cell(row, col) accessor (1-based)
# rather than with a single document-wide Find/Replace.
$rows = @(1, 5, 10, 15, 20)

$newValues = @(
    @("19×44=836", "20×70=1400", "22×75=1650", "91×18=1638", "25×42=1050"),
    @("60×21=1260", "59×81=4779", "69×42=2898", "34×61=2074", "11×72=792"),
    @("16×53=848", "44×74=3256", "82×46=3772", "25×21=525", "83×26=2158"),
    @("63×98=6174", "87×96=8352", "82×65=5330", "86×23=1978", "88×65=5720"),
    @("84×41=3444", "81×30=2430", "42×51=2142", "83×99=8217", "73×92=6716")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 1; $c -le 5; $c++) {
        $newText = $newValues[$r][$c - 1]
        $cellRange = $t.Cell($row, $c).Range
        # Exclude the trailing end-of-cell marker from the range we overwrite.
        $cellRange.MoveEnd(12, -1) | Out-Null
        $cellRange.Text = $newText
    }
}
